# Generate Report for Handoff
# Updates the localization-status report after a new handoff run:
#   - Overview sheet: refresh "Latest HO Xliff Generate Date" timestamps for the
#     e2e *.md rows that were just handed off (rows 7,9,11,12,13,14).
#   - zh-cn sheet: refresh "Latest Handoff Datetime" (col H) for the same rows,
#     and set their "Priority" (col E) to "ht" (handoff type) now that it has
#     been determined.
#   - de-de sheet: refresh "Latest Handoff Datetime" (col H) for the same rows,
#     and set their "Priority" (col E) to "ht" as well.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 11, 12, 13, 14)

foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-26 20:20:35"

    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-26 20:20:31"

    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-26 20:20:35"
}
